$d = $word.ActiveDocument

# Body paragraph: "This is gonna be a fancy word document with headers and footers"
#              -> "This is gonna bx a fancy word documxnt with hxadxrs and footxrs"
$body = $d.Paragraphs(1).Range
$body.Delete()
$body.InsertAfter("This is gonna bx a fancy word documxnt with hxadxrs and footxrs")

# Footer: "And im in the footer" -> "And im in thx footxr"
$ftr = $d.Sections(1).Footers(1).Range
$ftr.Delete()
$ftr.InsertAfter("And im in thx footxr")

# Header: "changing the header to this" -> "changing thx hxadxr to this"
$hdr = $d.Sections(1).Headers(1).Range
$hdr.Delete()
$hdr.InsertAfter("changing thx hxadxr to this")
